$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 10, shape "文本框 4": the three runs
#      "           " + "git" + " reset –hard  "
#    are merged back into a single run and the stray en-dash before "hard"
#    is corrected to a literal double-hyphen ("git reset --hard").
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(4)
$tr10 = $sh10.TextFrame.TextRange
$run10 = $tr10.Characters(14, 28)
$run10.Text = "           git reset --hard  "

# ---------------------------------------------------------------------------
# 2) Slide 11, shape "文本框 2": the run ' –C "Your Email"' is split into
#    two runs: ' –' and 'C "Your Email"' (ssh-keygen –t rsa –C ...).
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(1)
$tr11 = $sh11.TextFrame.TextRange
$dash11 = $tr11.Characters(53, 2)
$dash11.Font.Color.RGB = 255

# ---------------------------------------------------------------------------
# 3) Slide 6, shape "内容占位符 2": fix the description of `git config
#    --global` - it configures the *current* user, not *all* users.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange
$run6 = $tr6.Characters(90, 15)
$run6.Text = "配置对当前用户都普遍适用的配置"
